$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2102564102564103
$ws.Range("C2").Value = 0.5128205128205128
$ws.Range("J2").Value = 0.02051282051282051
$ws.Range("P2").Value = 0.1461538461538462
$ws.Range("S2").Value = 0.1102564102564103
$ws.Range("B3").Value = 0.009900990099009901
$ws.Range("C3").Value = 0.01485148514851485
$ws.Range("J3").Value = 0.03465346534653466
$ws.Range("O3").Value = 0.004950495049504951
$ws.Range("P3").Value = 0.7326732673267327
$ws.Range("S3").Value = 0.202970297029703
$ws.Range("J4").Value = 0.1020408163265306
$ws.Range("O4").Value = 0.02040816326530612
$ws.Range("P4").Value = 0.7551020408163265
$ws.Range("S4").Value = 0.1224489795918367
$ws.Range("B6").Value = 0.0945273631840796
$ws.Range("D6").Value = 0.01990049751243781
$ws.Range("E6").Value = 0.004975124378109453
$ws.Range("F6").Value = 0.0845771144278607
$ws.Range("J6").Value = 0.1890547263681592
$ws.Range("O6").Value = 0.03980099502487562
$ws.Range("Q6").Value = 0.1641791044776119
$ws.Range("R6").Value = 0.07960199004975124
$ws.Range("S6").Value = 0.3233830845771145
$ws.Range("B7").Value = 0.1290322580645161
$ws.Range("D7").Value = 0.02150537634408602
$ws.Range("F7").Value = 0.08602150537634409
$ws.Range("J7").Value = 0.1290322580645161
$ws.Range("O7").Value = 0.05376344086021505
$ws.Range("Q7").Value = 0.1344086021505376
$ws.Range("R7").Value = 0.05913978494623656
$ws.Range("S7").Value = 0.3870967741935484
$ws.Range("B8").Value = 0.1658878504672897
$ws.Range("D8").Value = 0.02803738317757009
$ws.Range("F8").Value = 0.05607476635514019
$ws.Range("J8").Value = 0.1004672897196262
$ws.Range("O8").Value = 0.0514018691588785
$ws.Range("Q8").Value = 0.1939252336448598
$ws.Range("R8").Value = 0.06542056074766354
$ws.Range("S8").Value = 0.338785046728972
$ws.Range("B9").Value = 0.1415094339622641
$ws.Range("D9").Value = 0.01886792452830189
$ws.Range("F9").Value = 0.07547169811320754
$ws.Range("J9").Value = 0.1037735849056604
$ws.Range("O9").Value = 0.05660377358490566
$ws.Range("Q9").Value = 0.2075471698113208
$ws.Range("R9").Value = 0.05660377358490566
$ws.Range("S9").Value = 0.3396226415094339
$ws.Range("B10").Value = 0.1381634372367312
$ws.Range("D10").Value = 0.02358887952822241
$ws.Range("F10").Value = 0.05897219882055602
$ws.Range("J10").Value = 0.117944397641112
$ws.Range("O10").Value = 0.03791069924178601
$ws.Range("Q10").Value = 0.2240943555181129
$ws.Range("R10").Value = 0.06908171861836562
$ws.Range("S10").Value = 0.3302443133951137
$ws.Range("G11").Value = 0.1254237288135593
$ws.Range("J11").Value = 0.0847457627118644
$ws.Range("K11").Value = 0.1796610169491525
$ws.Range("L11").Value = 0.6033898305084746
$ws.Range("S11").Value = 0.006779661016949152
$ws.Range("G12").Value = 0.7307692307692307
$ws.Range("J12").Value = 0.2087912087912088
$ws.Range("L12").Value = 0.02747252747252747
$ws.Range("S12").Value = 0.03296703296703297
$ws.Range("G13").Value = 0.6176470588235294
$ws.Range("J13").Value = 0.3529411764705883
$ws.Range("S13").Value = 0.02941176470588235
$ws.Range("F15").Value = 0.02164502164502164
$ws.Range("H15").Value = 0.1212121212121212
$ws.Range("I15").Value = 0.02597402597402598
$ws.Range("J15").Value = 0.3246753246753247
$ws.Range("K15").Value = 0.05194805194805195
$ws.Range("M15").Value = 0.008658008658008658
$ws.Range("O15").Value = 0.04761904761904762
$ws.Range("S15").Value = 0.3982683982683983
$ws.Range("F16").Value = 0.03361344537815126
$ws.Range("H16").Value = 0.1848739495798319
$ws.Range("I16").Value = 0.04201680672268908
$ws.Range("J16").Value = 0.4285714285714285
$ws.Range("K16").Value = 0.1134453781512605
$ws.Range("M16").Value = 0.008403361344537815
$ws.Range("O16").Value = 0.04621848739495799
$ws.Range("S16").Value = 0.1428571428571428
$ws.Range("F17").Value = 0.01415094339622642
$ws.Range("H17").Value = 0.160377358490566
$ws.Range("I17").Value = 0.04952830188679246
$ws.Range("J17").Value = 0.4882075471698113
$ws.Range("K17").Value = 0.09433962264150944
$ws.Range("M17").Value = 0.0259433962264151
$ws.Range("O17").Value = 0.06132075471698113
$ws.Range("S17").Value = 0.1061320754716981
$ws.Range("F18").Value = 0.02054794520547945
$ws.Range("H18").Value = 0.1780821917808219
$ws.Range("I18").Value = 0.07534246575342465
$ws.Range("J18").Value = 0.410958904109589
$ws.Range("K18").Value = 0.1506849315068493
$ws.Range("M18").Value = 0.00684931506849315
$ws.Range("O18").Value = 0.0410958904109589
$ws.Range("S18").Value = 0.1164383561643836
$ws.Range("F19").Value = 0.01633393829401089
$ws.Range("H19").Value = 0.2404718693284937
$ws.Range("I19").Value = 0.0499092558983666
$ws.Range("J19").Value = 0.3720508166969147
$ws.Range("K19").Value = 0.1252268602540835
$ws.Range("M19").Value = 0.0190562613430127
$ws.Range("O19").Value = 0.06261343012704174
$ws.Range("S19").Value = 0.1143375680580762
